$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Sun Jun 25 18:30:20 UTC 2023
# (includes a reorder: row 21 <-> row 22 swap between Uniswap and BinanceUSD)

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.469.88'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.73%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.890.62'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.14%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.37%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.94%  '

# Row 6
$ws.Range('E6').Value = '  +0.19%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4883'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2924'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.98%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06671'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.07%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.883.05'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.24%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '17.00'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.12%  '

# Row 12
$ws.Range('E12').Value = '  +1.59%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.130'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.71%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.27%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6646'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.01%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.432.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.74%  '

# Row 17
$ws.Range('E17').Value = '  +3.27%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007810'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.76%  '

# Row 19
$ws.Range('E19').Value = '  +0.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.125.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.14%  '

# Row 21
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.32%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.274'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +11.45%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '187.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.15%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.143'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.42%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.460'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.72%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.81%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.03%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.930'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.45%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.466'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.66%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.346'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.01%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09148'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.27%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.091'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.07%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05202'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.25%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7391'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.14%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.096'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.30%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.720'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.22%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01817'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.673'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.45%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9165'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.57%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.036'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.72%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4389'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.54%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.931'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.58%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '105.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.61%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9925'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.67%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1383'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.28%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '68.08'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +18.81%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.589'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.96%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.010'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.82%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.85'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.00%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05828'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.12%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3928'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.68%  '
